# Gaussian Quadrature Scheme export: rename sheet/tab to the short "CubeA"
# name, and append the new "HexGrid-60degTilt5degRes" scheme results as
# row 16 of the averaged-intensities table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the worksheet (tab) --------------------------------------
$ws.Name = "CubeA"

# --- Append the new data row (row 16) ---------------------------------
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 0.8620159336131725
$ws.Cells.Item(16, 4).Value = 1.1851820878893
$ws.Cells.Item(16, 5).Value = 0.9528535474453399
$ws.Cells.Item(16, 6).Value = 1.036449921668027
$ws.Cells.Item(16, 7).Value = 0.8620159336131725
$ws.Cells.Item(16, 8).Value = 1.1851820878893
$ws.Cells.Item(16, 9).Value = 0.9329395222185278
$ws.Cells.Item(16, 10).Value = 1.033450248940469
$ws.Cells.Item(16, 11).Value = 0.9489456796225956
$ws.Cells.Item(16, 12).Value = 1.11605938060394
$ws.Cells.Item(16, 13).Value = 0.8620159336131725
$ws.Cells.Item(16, 14).Value = 1.06901781766732
$ws.Cells.Item(16, 15).Value = 1.00912537265396
$ws.Cells.Item(16, 16).Value = 1.008487040250172

# Match the existing formatting of column A (bold/centered/bordered header
# style) by copying it from the row above, same as every other row in the
# table.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Tiny floating point tweaks to three pre-existing cells -----------
$ws.Cells.Item(15, 4).Value = 0.9096925873334168
$ws.Cells.Item(15, 8).Value = 0.9096925873334168
$ws.Cells.Item(15, 15).Value = 0.9887980053438663
